# Add files via upload
# - Clear stale "実施内容" (F-column) entries for rows 18-22 and 25-29
# - Record progress for the "Nodejs MongoDb Express" course (row 13) and add
#   two new "React系列" course entries for rows 14 and 15, each with a 100%
#   progress mark in column G
# - Move the active cell selection to H16

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear out the old plan entries that no longer apply
$ws.Range("F18").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("F22").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("F29").Value = ""

# Row 13: mark the Nodejs MongoDb Express 34-39 entry as completed
$ws.Range("G13").Value = 1
$ws.Range("G13").NumberFormat = "0%"

# Row 14: new React course entry, marked completed
$ws.Range("F14").Value = "React系列课程从零基础到项目开发实战 1-19"
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = "0%"

# Row 15: new React course entry, marked completed
$ws.Range("F15").Value = "React系列课程从零基础到项目开发实战 20-36"
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = "0%"

# Update the saved selection to H16
$ws.Range("H16").Select() | Out-Null
